# Simulated Wild Card round and logged it
# Update "R" (road) row depth-of-target stats on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 ("R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 296
$wsOff.Range("C3").Value = 197
$wsOff.Range("D3").Value = 105
$wsOff.Range("E3").Value = 45
$wsOff.Range("F3").Value = 5
$wsOff.Range("G3").Value = 7

# DEF sheet - row 3 ("R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 533
$wsDef.Range("C3").Value = 399
$wsDef.Range("D3").Value = 82
$wsDef.Range("E3").Value = 40
